$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values per repulled data
$ws.Range("F4").Value = -4
$ws.Range("F5").Value = 2
$ws.Range("F7").Value = -1
$ws.Range("F8").Value = -7
$ws.Range("F10").Value = -4
$ws.Range("F11").Value = 2
$ws.Range("F15").Value = 4
$ws.Range("F19").Value = -2
$ws.Range("F20").Value = -5
$ws.Range("F21").Value = -8
$ws.Range("F24").Value = -4
$ws.Range("F26").Value = 0
$ws.Range("F29").Value = 1
